# Přidána sekce "Společný HW - detailní popis IO atd" -> nový list commonHW_DI
# s parametry digitálních vstupů (DI), umístěný za poslední list "TGZpMotion".

$wb = $excel.ActiveWorkbook

# --- 1. Přidat nový list na konec sešitu a pojmenovat jej -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "commonHW_DI"

# --- 2. Šířky sloupců (OOXML "width" = ColumnWidth + 5/6) -------------------
$colWidths = @{ 1=10.2; 2=20.35; 3=20.35; 4=23.27; 5=22.02; 6=26.2; 7=26.2; 8=31.62; 9=27.17; 10=15.34; 11=30.23 }
foreach ($c in $colWidths.Keys) {
    $ws.Columns.Item($c).ColumnWidth = $colWidths[$c] - (5.0/6.0)
}

# --- 3. Hlavička (řádky 1-3) -------------------------------------------------
$ws.Cells.Item(1,1).Value  = " č. vstupu "
$ws.Cells.Item(1,2).Value  = " rozsah napětí log. 0 "
$ws.Cells.Item(1,3).Value  = " rozsah napětí log. 1 "
$ws.Cells.Item(1,4).Value  = " nominální vstupní napětí "
$ws.Cells.Item(1,5).Value  = " rozsah napájení DI/DO "
$ws.Cells.Item(1,6).Value  = " spotřeba v log. 1 (17V) ±20% "
$ws.Cells.Item(1,7).Value  = " spotřeba v log. 1 (24V) ±20% "
$ws.Cells.Item(1,8).Value  = " vlastní spotřeba v log.1 (28V) ±20% "
$ws.Cells.Item(1,9).Value  = " Nominální vstupní odpor ±20%"
$ws.Cells.Item(1,10).Value = " Přiřazeno k ose "
$ws.Cells.Item(1,11).Value = " Max. vstupní frekvence - obdélník "

$ws.Cells.Item(2,1).Value  = " # "
$ws.Cells.Item(2,2).Value  = " U<sub>log0</sub> "
$ws.Cells.Item(2,3).Value  = " U<sub>log1</sub> "
$ws.Cells.Item(2,4).Value  = " U<sub>nom</sub>"
$ws.Cells.Item(2,5).Value  = " VDD<sub>IO</sub>"
$ws.Cells.Item(2,6).Value  = " I<sub>in17</sub>"
$ws.Cells.Item(2,7).Value  = " I<sub>in24</sub>"
$ws.Cells.Item(2,8).Value  = " I<sub>in28</sub>"
$ws.Cells.Item(2,9).Value  = " RI<sub>inNom</sub> "
$ws.Cells.Item(2,10).Value = " Osa č. "
$ws.Cells.Item(2,11).Value = " f<sub>maxSq</sub>"

$ws.Cells.Item(3,1).Value  = "-"
$ws.Cells.Item(3,2).Value  = " V"
$ws.Cells.Item(3,3).Value  = " V"
$ws.Cells.Item(3,4).Value  = " V"
$ws.Cells.Item(3,5).Value  = " V "
$ws.Cells.Item(3,6).Value  = " mA"
$ws.Cells.Item(3,7).Value  = " mA"
$ws.Cells.Item(3,8).Value  = " mA"
$ws.Cells.Item(3,9).Value  = " kΩ "
$ws.Cells.Item(3,10).Value = "-`t`t`t`t`t "
$ws.Cells.Item(3,11).Value = "kHz`t`t`t`t`t "

# --- 4. Datové řádky 4-11 (8 digitálních vstupů) -----------------------------
# Společné sloupce B/C/D jsou pro všechny vstupy stejné.
# Vstupy 1-6 (řádky 4-9) mají jinou spotřebu/odpor než vstupy 7-8 (řádky 10-11,
# ty navíc nevyžadují přiřazení osy).
$rows = @(
    @{ r=4;  a=1; e=" 17-28V"; f=2.05; g=2.92; h=4.15; i=8.2; j=1 },
    @{ r=5;  a=2; e=" 17-28V"; f=2.05; g=2.92; h=4.15; i=8.2; j=2 },
    @{ r=6;  a=3; e=" 17-28V"; f=2.05; g=2.92; h=4.15; i=8.2; j=1 },
    @{ r=7;  a=4; e=" 17-28V"; f=2.05; g=2.92; h=4.15; i=8.2; j=2 },
    @{ r=8;  a=5; e=" 17-28V"; f=2.05; g=2.92; h=4.15; i=8.2; j=1 },
    @{ r=9;  a=6; e=" 17-28V"; f=2.05; g=2.92; h=4.15; i=8.2; j=2 },
    @{ r=10; a=7; e=" nevyžaduje"; f=1.25; g=3.32; h=4.5; i=7.2; j=1 },
    @{ r=11; a=8; e=" nevyžaduje"; f=1.25; g=3.32; h=4.5; i=7.2; j=2 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r,1).Value  = $row.a
    $ws.Cells.Item($r,2).Value  = " 0-10 "
    $ws.Cells.Item($r,3).Value  = " 17-28V "
    $ws.Cells.Item($r,4).Value  = 24
    $ws.Cells.Item($r,5).Value  = $row.e
    $ws.Cells.Item($r,6).Value  = $row.f
    $ws.Cells.Item($r,7).Value  = $row.g
    $ws.Cells.Item($r,8).Value  = $row.h
    $ws.Cells.Item($r,9).Value  = $row.i
    $ws.Cells.Item($r,10).Value = $row.j
    $ws.Cells.Item($r,11).Value = 50
}

# --- 5. Prázdné naformátované řádky 12-21 ve sloupcích A:B (zachování rozsahu
#        použité oblasti A1:K21 jako ve zdrojovém souboru) ---------------------
$ws.Range("A12:B21").Font.Name = "Arial"

# --- 6. Aktivní buňka / výběr na novém listu (dle zdrojového souboru G20) ----
$ws.Range("G20").Select()
